# Remove the stray/orphan bookmark "_k7637xdazqek" from the document.
# Word automatically renumbers the w:id of every remaining bookmark
# (bookmarkStart/bookmarkEnd pairs) downward to close the gap, which is
# exactly what the target revision shows: the bookmark is gone and all
# following bookmarks shift id N -> N-1 while keeping their names.
$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_k7637xdazqek")) {
    $d.Bookmarks.Item("_k7637xdazqek").Delete()
}
